$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.353290908
$ws.Range("C2").Value = 5.104107506
$ws.Range("D2").Value = 4.570819855
$ws.Range("E2").Value = 2.758747578
$ws.Range("F2").Value = 0.7508165980000001
$ws.Range("G2").Value = 0.5049157875049435
$ws.Range("H2").Value = 0.0803598434265488
$ws.Range("I2").Value = 8.368468843012426

$ws.Range("B3").Value = 7.054631072
$ws.Range("C3").Value = 7.804447965
$ws.Range("D3").Value = 0.2230691910000001
$ws.Range("E3").Value = 0.1952342987
$ws.Range("F3").Value = 0.7498168929999993
$ws.Range("G3").Value = 0.133281627780859
$ws.Range("H3").Value = 0.02121242988465779
$ws.Range("I3").Value = 8.379626233867196

$ws.Range("B4").Value = 9.004154965
$ws.Range("C4").Value = 9.753971722999999
$ws.Range("D4").Value = 0.04548645019999999
$ws.Range("E4").Value = 0.0364074707
$ws.Range("F4").Value = 0.7498167579999997
$ws.Range("G4").Value = 0.2226404912090963
$ws.Range("H4").Value = 0.03543433470833535
$ws.Range("I4").Value = 8.379627742568523
